$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue")
$ws.Range("B16").Value = 12
